$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Top_YTD")

$data = @(
    @("SITAB CI (STBC)", 168.36),
    @("UNILEVER CI (UNLC)", 126.18),
    @("FILTISAC CI (FTSC)", 119.31),
    @("TRACTAFRIC MOTORS CI (PRSC)", 105.01),
    @("BANK OF AFRICA SENEGAL (BOAS)", 63.65),
    @("SAPH CI (SPHC)", 62.49),
    @("SOCIETE IVOIRIENNE DE BANQUE  (SIBC)", 59.76),
    @("TOTALENERGIES MARKETING CI (TTLC)", 48.3),
    @("UNIWAX CI (UNXC)", 45.79),
    @("BICI CI (BICC)", 38.88)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $row++
}

$wb.Save()
